# "change shortest path length; and update result"
#
# Row 8 / Row 9 of the "Network characteristic" / "Top3" table get new
# labels + new result lists:
#   B8: "average shortest length" -> "co-conservation"
#   C8: "253D,190R,493Q"          -> "614D-1027T,614D-859T,376T-408R"
#   B9: "co-conservation"         -> "shortest path length"
#   C9: "614D-1027T,614D-859T,376T-408R" -> "190R-679N,5L-253D,190R-484E"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B8").Value = "co-conservation"
$ws.Range("C8").Value = "614D-1027T,614D-859T,376T-408R"

$ws.Range("B9").Value = "shortest path length"
$ws.Range("C9").Value = "190R-679N,5L-253D,190R-484E"
